$wb = $excel.ActiveWorkbook

# --- Update the time_taken timestamps on the "data" sheet (column F, rows 2-13) ---
$ds = $wb.Worksheets.Item("data")

$ds.Range("F2").Value  = "2021-10-05 14:35:06.215092"
$ds.Range("F3").Value  = "2021-10-05 14:35:06.215100"
$ds.Range("F4").Value  = "2021-10-05 14:35:06.215103"
$ds.Range("F5").Value  = "2021-10-05 14:35:06.215106"
$ds.Range("F6").Value  = "2021-10-05 14:35:06.215109"
$ds.Range("F7").Value  = "2021-10-05 14:35:06.215112"
$ds.Range("F8").Value  = "2021-10-05 14:35:06.215115"
$ds.Range("F9").Value  = "2021-10-05 14:35:06.215118"
$ds.Range("F10").Value = "2021-10-05 14:35:06.215120"
$ds.Range("F11").Value = "2021-10-05 14:35:06.215123"
$ds.Range("F12").Value = "2021-10-05 14:35:06.215126"
$ds.Range("F13").Value = "2021-10-05 14:35:06.215128"

# --- Add a new "metadata" worksheet placed right after "data" ---
$ws = $wb.Worksheets.Add($null, $ds)
$ws.Name = "metadata"

# Match the outline/summary layout used by the "data" sheet
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Match the page margins used by the "data" sheet (values are in points: 72pt = 1 inch)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row (B1:G1) -- bold, thin border all around, centered horizontally & top-aligned
# (matches the look of the "data" sheet's own header row / index column)
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# A2 (row index) -- same styling as the header
$ws.Range("A2").Value = 0
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160

# Data row
$ws.Range("B2").Value = "Oligodontia"
$ws.Range("C2").Value = 148

# Leading apostrophe forces this to be stored as text "0.6" rather than
# being interpreted/coerced as the number 0.6
$ws.Range("D2").Value = "'0.6"

$ws.Range("E2").Value = "2021-09-13T10:24:10.790328Z"
$ws.Range("F2").Value = "2021-10-05 14:35:06.211279"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/148/?format=json"

$ds.Select()
$ds.Range("A1").Select()
